# "time record + record cheat"
#  1) The auto-updating "datetimeFigureOut" date field cached on the
#     slide master and on every slide layout advances by one Hebrew day:
#     ט"ז/ניסן/תשפ"ב (16 Nisan 5782)  ->  י"ז/ניסן/תשפ"ב (17 Nisan 5782)
#  2) The big decorative oval on slide 1 gets a darker purple fill:
#     A47DC5 -> 805C9C

$p = $ppt.ActivePresentation

$oldDate = 'ט"ז/ניסן/תשפ"ב'
$newDate = 'י"ז/ניסן/תשפ"ב'

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        $isDate = $false
        try {
            if ($sh.HasTextFrame -and $sh.PlaceholderFormat.Type -eq 16) {
                $isDate = $true
            }
        } catch {
            $isDate = $false
        }
        if (-not $isDate) {
            if ($sh.Name -like 'Date Placeholder*') {
                $isDate = $true
            }
        }
        if ($isDate) {
            $tr = $sh.TextFrame.TextRange
            if ($tr.Text -eq $oldDate) {
                $tr.Text = $newDate
            }
        }
    }
}

# 1a) Slide master
Update-DatePlaceholder $p.SlideMaster.Shapes

# 1b) Every slide layout hanging off the master
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DatePlaceholder $layouts.Item($li).Shapes
}

# 2) Recolor the big oval (fill) on slide 1
$slide1 = $p.Slides.Item(1)
for ($i = 1; $i -le $slide1.Shapes.Count; $i++) {
    $sh = $slide1.Shapes.Item($i)
    if ($sh.Name -eq 'Oval 8') {
        $sh.Fill.ForeColor.RGB = 10247296   # 0x9C5C80 little-endian == RGB(805C9C)
    }
}
